{"js": "// Replace each arithmetic expression in the single 20x5 table with its\n// updated value, one cell at a time. The mapping below records the\n// (row, col) position together with the expected old text (for a sanity\n// check) and the new text that should replace it, in document order --\n// this mirrors the diff, which only rewrites the <w:t> text runs inside\n// each table cell and leaves every other property (fonts, size,\n// paragraph alignment, etc.) untouched.\nconst replacements = [\n  [0, 0, \"9+80=89\", \"16+34=50\"],\n  [0, 1, \"1+70=71\", \"79-27=52\"],\n  [0, 2, \"10+75=85\", \"90-21=69\"],\n  [0, 3, \"95-44=51\", \"71-13=58\"],\n  [0, 4, \"5+7=12\", \"87+12=99\"],\n  [1, 0, \"11+81=92\", \"44+39=83\"],\n  [1, 1, \"94-64=30\", \"99-2=97\"],\n  [1, 2, \"65-8=57\", \"52+43=95\"],\n  [1, 3, \"47+25=72\", \"30+34=64\"],\n  [1, 4, \"49-40=9\", \"32+14=46\"],\n  [2, 0, \"26+28=54\", \"56-39=17\"],\n  [2, 1, \"90+2=92\", \"40-3=37\"],\n  [2, 2, \"51-18=33\", \"7+13=20\"],\n  [2, 3, \"68+4=72\", \"55+4=59\"],\n  [2, 4, \"89-4=85\", \"36-3=33\"],\n  [3, 0, \"10+61=71\", \"39+52=91\"],\n  [3, 1, \"91-77=14\", \"18+49=67\"],\n  [3, 2, \"73-16=57\", \"60-15=45\"],\n  [3, 3, \"93-0=93\", \"88-24=64\"],\n  [3, 4, \"46+36=82\", \"69-7=62\"],\n  [4, 0, \"18+55=73\", \"12+13=25\"],\n  [4, 1, \"34+21=55\", \"10+35=45\"],\n  [4, 2, \"86-69=17\", \"84-16=68\"],\n  [4, 3, \"88-26=62\", \"22+6=28\"],\n  [4, 4, \"36+3=39\", \"39-34=5\"],\n  [5, 0, \"74+22=96\", \"31+17=48\"],\n  [5, 1, \"38-37=1\", \"38+42=80\"],\n  [5, 2, \"82-45=37\", \"61-22=39\"],\n  [5, 3, \"0+50=50\", \"86+10=96\"],\n  [5, 4, \"13+8=21\", \"52-41=11\"],\n  [6, 0, \"47-20=27\", \"79-38=41\"],\n  [6, 1, \"49-43=6\", \"52-25=27\"],\n  [6, 2, \"33+51=84\", \"66-18=48\"],\n  [6, 3, \"69-12=57\", \"0+49=49\"],\n  [6, 4, \"87-80=7\", \"18+59=77\"],\n  [7, 0, \"34+46=80\", \"69+16=85\"],\n  [7, 1, \"72-55=17\", \"28-14=14\"],\n  [7, 2, \"71-38=33\", \"99-35=64\"],\n  [7, 3, \"76-21=55\", \"32+52=84\"],\n  [7, 4, \"30+15=45\", \"31+2=33\"],\n  [8, 0, \"19-5=14\", \"73-52=21\"],\n  [8, 1, \"29+65=94\", \"46+30=76\"],\n  [8, 2, \"71+19=90\", \"67-59=8\"],\n  [8, 3, \"2+93=95\", \"2+85=87\"],\n  [8, 4, \"39+18=57\", \"4+34=38\"],\n  [9, 0, \"78-54=24\", \"23-10=13\"],\n  [9, 1, \"4+24=28\", \"88-43=45\"],\n  [9, 2, \"65+24=89\", \"30+46=76\"],\n  [9, 3, \"35+47=82\", \"21+13=34\"],\n  [9, 4, \"30-18=12\", \"55+27=82\"],\n  [10, 0, \"61+35=96\", \"35+48=83\"],\n  [10, 1, \"34+19=53\", \"95-86=9\"],\n  [10, 2, \"47-46=1\", \"5+73=78\"],\n  [10, 3, \"82-1=81\", \"58-4=54\"],\n  [10, 4, \"2+32=34\", \"87-61=26\"],\n  [11, 0, \"1+87=88\", \"37-8=29\"],\n  [11, 1, \"5+11=16\", \"83-4=79\"],\n  [11, 2, \"41-0=41\", \"67+15=82\"],\n  [11, 3, \"66-44=22\", \"63-34=29\"],\n  [11, 4, \"89-87=2\", \"85-47=38\"],\n  [12, 0, \"36+41=77\", \"15+50=65\"],\n  [12, 1, \"24+36=60\", \"83-22=61\"],\n  [12, 2, \"94-48=46\", \"14+53=67\"],\n  [12, 3, \"33+19=52\", \"7+69=76\"],\n  [12, 4, \"12+13=25\", \"20+39=59\"],\n  [13, 0, \"21+2=23\", \"62-12=50\"],\n  [13, 1, \"79-40=39\", \"80-62=18\"],\n  [13, 2, \"94-46=48\", \"40-34=6\"],\n  [13, 3, \"76-56=20\", \"49-18=31\"],\n  [13, 4, \"97-23=74\", \"30+19=49\"],\n  [14, 0, \"10+29=39\", \"25-21=4\"],\n  [14, 1, \"30-0=30\", \"41-26=15\"],\n  [14, 2, \"97-43=54\", \"81-39=42\"],\n  [14, 3, \"39+37=76\", \"10+1=11\"],\n  [14, 4, \"6+17=23\", \"56+3=59\"],\n  [15, 0, \"76-72=4\", \"12-11=1\"],\n  [15, 1, \"60-20=40\", \"70-67=3\"],\n  [15, 2, \"43-18=25\", \"51-33=18\"],\n  [15, 3, \"63-3=60\", \"91-6=85\"],\n  [15, 4, \"10+17=27\", \"34+12=46\"],\n  [16, 0, \"23+67=90\", \"23+69=92\"],\n  [16, 1, \"41-12=29\", \"48-41=7\"],\n  [16, 2, \"51+19=70\", \"20-9=11\"],\n  [16, 3, \"32+37=69\", \"92-54=38\"],\n  [16, 4, \"90-61=29\", \"31-17=14\"],\n  [17, 0, \"44+44=88\", \"6+69=75\"],\n  [17, 1, \"15+33=48\", \"66+6=72\"],\n  [17, 2, \"50+5=55\", \"46+9=55\"],\n  [17, 3, \"35+45=80\", \"68+23=91\"],\n  [17, 4, \"99-52=47\", \"55-30=25\"],\n  [18, 0, \"24+25=49\", \"40+58=98\"],\n  [18, 1, \"20+33=53\", \"69-15=54\"],\n  [18, 2, \"52+19=71\", \"50+33=83\"],\n  [18, 3, \"67+1=68\", \"90-12=78\"],\n  [18, 4, \"3+68=71\", \"3+79=82\"],\n  [19, 0, \"69-1=68\", \"16+3=19\"],\n  [19, 1, \"69-55=14\", \"71+16=87\"],\n  [19, 2, \"8+24=32\", \"36+21=57\"],\n  [19, 3, \"15+83=98\", \"13+78=91\"],\n  [19, 4, \"8+27=35\", \"88-47=41\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\n\nconst table = tables.items[0];\n\n// Resolve every cell's first paragraph up front so we can batch the\n// property loads before mutating anything.\nconst paragraphs = [];\nfor (const [row, col] of replacements) {\n  const cell = table.getCell(row, col);\n  const paragraph = cell.body.paragraphs.getFirst();\n  paragraph.load(\"text\");\n  paragraphs.push(paragraph);\n}\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [row, col, oldText, newText] = replacements[i];\n  const paragraph = paragraphs[i];\n  // Defensive check -- only replace if the cell still holds the text we\n  // expect; otherwise leave it untouched rather than risk clobbering the\n  // wrong cell.\n  if (paragraph.text === oldText) {\n    const range = paragraph.getRange();\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each arithmetic expression in the single 20x5 table with its\n# updated value, one cell at a time. The mapping below records the\n# (Row, Col) position together with the expected old text (for a sanity\n# check) and the new text that should replace it, in document order --\n# this mirrors the diff, which only rewrites the text inside each table\n# cell and leaves every other property (fonts, size, paragraph\n# alignment, etc.) untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Row = 1; Col = 1; OldText = \"9+80=89\"; NewText = \"16+34=50\" }\n    @{ Row = 1; Col = 2; OldText = \"1+70=71\"; NewText = \"79-27=52\" }\n    @{ Row = 1; Col = 3; OldText = \"10+75=85\"; NewText = \"90-21=69\" }\n    @{ Row = 1; Col = 4; OldText = \"95-44=51\"; NewText = \"71-13=58\" }\n    @{ Row = 1; Col = 5; OldText = \"5+7=12\"; NewText = \"87+12=99\" }\n    @{ Row = 2; Col = 1; OldText = \"11+81=92\"; NewText = \"44+39=83\" }\n    @{ Row = 2; Col = 2; OldText = \"94-64=30\"; NewText = \"99-2=97\" }\n    @{ Row = 2; Col = 3; OldText = \"65-8=57\"; NewText = \"52+43=95\" }\n    @{ Row = 2; Col = 4; OldText = \"47+25=72\"; NewText = \"30+34=64\" }\n    @{ Row = 2; Col = 5; OldText = \"49-40=9\"; NewText = \"32+14=46\" }\n    @{ Row = 3; Col = 1; OldText = \"26+28=54\"; NewText = \"56-39=17\" }\n    @{ Row = 3; Col = 2; OldText = \"90+2=92\"; NewText = \"40-3=37\" }\n    @{ Row = 3; Col = 3; OldText = \"51-18=33\"; NewText = \"7+13=20\" }\n    @{ Row = 3; Col = 4; OldText = \"68+4=72\"; NewText = \"55+4=59\" }\n    @{ Row = 3; Col = 5; OldText = \"89-4=85\"; NewText = \"36-3=33\" }\n    @{ Row = 4; Col = 1; OldText = \"10+61=71\"; NewText = \"39+52=91\" }\n    @{ Row = 4; Col = 2; OldText = \"91-77=14\"; NewText = \"18+49=67\" }\n    @{ Row = 4; Col = 3; OldText = \"73-16=57\"; NewText = \"60-15=45\" }\n    @{ Row = 4; Col = 4; OldText = \"93-0=93\"; NewText = \"88-24=64\" }\n    @{ Row = 4; Col = 5; OldText = \"46+36=82\"; NewText = \"69-7=62\" }\n    @{ Row = 5; Col = 1; OldText = \"18+55=73\"; NewText = \"12+13=25\" }\n    @{ Row = 5; Col = 2; OldText = \"34+21=55\"; NewText = \"10+35=45\" }\n    @{ Row = 5; Col = 3; OldText = \"86-69=17\"; NewText = \"84-16=68\" }\n    @{ Row = 5; Col = 4; OldText = \"88-26=62\"; NewText = \"22+6=28\" }\n    @{ Row = 5; Col = 5; OldText = \"36+3=39\"; NewText = \"39-34=5\" }\n    @{ Row = 6; Col = 1; OldText = \"74+22=96\"; NewText = \"31+17=48\" }\n    @{ Row = 6; Col = 2; OldText = \"38-37=1\"; NewText = \"38+42=80\" }\n    @{ Row = 6; Col = 3; OldText = \"82-45=37\"; NewText = \"61-22=39\" }\n    @{ Row = 6; Col = 4; OldText = \"0+50=50\"; NewText = \"86+10=96\" }\n    @{ Row = 6; Col = 5; OldText = \"13+8=21\"; NewText = \"52-41=11\" }\n    @{ Row = 7; Col = 1; OldText = \"47-20=27\"; NewText = \"79-38=41\" }\n    @{ Row = 7; Col = 2; OldText = \"49-43=6\"; NewText = \"52-25=27\" }\n    @{ Row = 7; Col = 3; OldText = \"33+51=84\"; NewText = \"66-18=48\" }\n    @{ Row = 7; Col = 4; OldText = \"69-12=57\"; NewText = \"0+49=49\" }\n    @{ Row = 7; Col = 5; OldText = \"87-80=7\"; NewText = \"18+59=77\" }\n    @{ Row = 8; Col = 1; OldText = \"34+46=80\"; NewText = \"69+16=85\" }\n    @{ Row = 8; Col = 2; OldText = \"72-55=17\"; NewText = \"28-14=14\" }\n    @{ Row = 8; Col = 3; OldText = \"71-38=33\"; NewText = \"99-35=64\" }\n    @{ Row = 8; Col = 4; OldText = \"76-21=55\"; NewText = \"32+52=84\" }\n    @{ Row = 8; Col = 5; OldText = \"30+15=45\"; NewText = \"31+2=33\" }\n    @{ Row = 9; Col = 1; OldText = \"19-5=14\"; NewText = \"73-52=21\" }\n    @{ Row = 9; Col = 2; OldText = \"29+65=94\"; NewText = \"46+30=76\" }\n    @{ Row = 9; Col = 3; OldText = \"71+19=90\"; NewText = \"67-59=8\" }\n    @{ Row = 9; Col = 4; OldText = \"2+93=95\"; NewText = \"2+85=87\" }\n    @{ Row = 9; Col = 5; OldText = \"39+18=57\"; NewText = \"4+34=38\" }\n    @{ Row = 10; Col = 1; OldText = \"78-54=24\"; NewText = \"23-10=13\" }\n    @{ Row = 10; Col = 2; OldText = \"4+24=28\"; NewText = \"88-43=45\" }\n    @{ Row = 10; Col = 3; OldText = \"65+24=89\"; NewText = \"30+46=76\" }\n    @{ Row = 10; Col = 4; OldText = \"35+47=82\"; NewText = \"21+13=34\" }\n    @{ Row = 10; Col = 5; OldText = \"30-18=12\"; NewText = \"55+27=82\" }\n    @{ Row = 11; Col = 1; OldText = \"61+35=96\"; NewText = \"35+48=83\" }\n    @{ Row = 11; Col = 2; OldText = \"34+19=53\"; NewText = \"95-86=9\" }\n    @{ Row = 11; Col = 3; OldText = \"47-46=1\"; NewText = \"5+73=78\" }\n    @{ Row = 11; Col = 4; OldText = \"82-1=81\"; NewText = \"58-4=54\" }\n    @{ Row = 11; Col = 5; OldText = \"2+32=34\"; NewText = \"87-61=26\" }\n    @{ Row = 12; Col = 1; OldText = \"1+87=88\"; NewText = \"37-8=29\" }\n    @{ Row = 12; Col = 2; OldText = \"5+11=16\"; NewText = \"83-4=79\" }\n    @{ Row = 12; Col = 3; OldText = \"41-0=41\"; NewText = \"67+15=82\" }\n    @{ Row = 12; Col = 4; OldText = \"66-44=22\"; NewText = \"63-34=29\" }\n    @{ Row = 12; Col = 5; OldText = \"89-87=2\"; NewText = \"85-47=38\" }\n    @{ Row = 13; Col = 1; OldText = \"36+41=77\"; NewText = \"15+50=65\" }\n    @{ Row = 13; Col = 2; OldText = \"24+36=60\"; NewText = \"83-22=61\" }\n    @{ Row = 13; Col = 3; OldText = \"94-48=46\"; NewText = \"14+53=67\" }\n    @{ Row = 13; Col = 4; OldText = \"33+19=52\"; NewText = \"7+69=76\" }\n    @{ Row = 13; Col = 5; OldText = \"12+13=25\"; NewText = \"20+39=59\" }\n    @{ Row = 14; Col = 1; OldText = \"21+2=23\"; NewText = \"62-12=50\" }\n    @{ Row = 14; Col = 2; OldText = \"79-40=39\"; NewText = \"80-62=18\" }\n    @{ Row = 14; Col = 3; OldText = \"94-46=48\"; NewText = \"40-34=6\" }\n    @{ Row = 14; Col = 4; OldText = \"76-56=20\"; NewText = \"49-18=31\" }\n    @{ Row = 14; Col = 5; OldText = \"97-23=74\"; NewText = \"30+19=49\" }\n    @{ Row = 15; Col = 1; OldText = \"10+29=39\"; NewText = \"25-21=4\" }\n    @{ Row = 15; Col = 2; OldText = \"30-0=30\"; NewText = \"41-26=15\" }\n    @{ Row = 15; Col = 3; OldText = \"97-43=54\"; NewText = \"81-39=42\" }\n    @{ Row = 15; Col = 4; OldText = \"39+37=76\"; NewText = \"10+1=11\" }\n    @{ Row = 15; Col = 5; OldText = \"6+17=23\"; NewText = \"56+3=59\" }\n    @{ Row = 16; Col = 1; OldText = \"76-72=4\"; NewText = \"12-11=1\" }\n    @{ Row = 16; Col = 2; OldText = \"60-20=40\"; NewText = \"70-67=3\" }\n    @{ Row = 16; Col = 3; OldText = \"43-18=25\"; NewText = \"51-33=18\" }\n    @{ Row = 16; Col = 4; OldText = \"63-3=60\"; NewText = \"91-6=85\" }\n    @{ Row = 16; Col = 5; OldText = \"10+17=27\"; NewText = \"34+12=46\" }\n    @{ Row = 17; Col = 1; OldText = \"23+67=90\"; NewText = \"23+69=92\" }\n    @{ Row = 17; Col = 2; OldText = \"41-12=29\"; NewText = \"48-41=7\" }\n    @{ Row = 17; Col = 3; OldText = \"51+19=70\"; NewText = \"20-9=11\" }\n    @{ Row = 17; Col = 4; OldText = \"32+37=69\"; NewText = \"92-54=38\" }\n    @{ Row = 17; Col = 5; OldText = \"90-61=29\"; NewText = \"31-17=14\" }\n    @{ Row = 18; Col = 1; OldText = \"44+44=88\"; NewText = \"6+69=75\" }\n    @{ Row = 18; Col = 2; OldText = \"15+33=48\"; NewText = \"66+6=72\" }\n    @{ Row = 18; Col = 3; OldText = \"50+5=55\"; NewText = \"46+9=55\" }\n    @{ Row = 18; Col = 4; OldText = \"35+45=80\"; NewText = \"68+23=91\" }\n    @{ Row = 18; Col = 5; OldText = \"99-52=47\"; NewText = \"55-30=25\" }\n    @{ Row = 19; Col = 1; OldText = \"24+25=49\"; NewText = \"40+58=98\" }\n    @{ Row = 19; Col = 2; OldText = \"20+33=53\"; NewText = \"69-15=54\" }\n    @{ Row = 19; Col = 3; OldText = \"52+19=71\"; NewText = \"50+33=83\" }\n    @{ Row = 19; Col = 4; OldText = \"67+1=68\"; NewText = \"90-12=78\" }\n    @{ Row = 19; Col = 5; OldText = \"3+68=71\"; NewText = \"3+79=82\" }\n    @{ Row = 20; Col = 1; OldText = \"69-1=68\"; NewText = \"16+3=19\" }\n    @{ Row = 20; Col = 2; OldText = \"69-55=14\"; NewText = \"71+16=87\" }\n    @{ Row = 20; Col = 3; OldText = \"8+24=32\"; NewText = \"36+21=57\" }\n    @{ Row = 20; Col = 4; OldText = \"15+83=98\"; NewText = \"13+78=91\" }\n    @{ Row = 20; Col = 5; OldText = \"8+27=35\"; NewText = \"88-47=41\" }\n)\n\n$tbl = $d.Tables.Item(1)\n\nforeach ($r in $replacements) {\n    $cell = $tbl.Cell($r.Row, $r.Col)\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    # Defensive check -- only replace if the cell still holds the text we\n    # expect; otherwise leave it untouched rather than risk clobbering the\n    # wrong cell.\n    if ($current -eq $r.OldText) {\n        $cell.Range.Text = $r.NewText\n    }\n}\n"}
